$d = $word.ActiveDocument

# 1) Update the first bullet text (Controls -> Hold WASD ... )
$null = $d.Content.Find.Execute('Hold W,A,S, and D to move up, left, down, and right respectively', $true, $false, $false, $false, $false, $true, 1, $false, 'Hold W,A,S, and D to move up, left, down, and right respectively, also should see animation playing.', 2)

# 2) Fill in the previously-empty last paragraph with the "Press T" bullet,
#    then append the rest of the new test-plan bullets after it.
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Press T -> tutorial pop up on the screen in loaded font'
# first item is ilvl 0 already (inherited from the template paragraph)

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Press G -> debug mode is enabled. The user should see all aabb colliders of entities being shown as red rectangle.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Move cursor around -> player character should always look at the cursor no matter of their walking direction.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Stay still and wait for enemy spawn -> wolves should shoot three bullets at a time (shotgun), bees should shoot one bullet at a time, and bombers should try to hit the player with their body. Bullets cause one damage and bombers cause two. They will be deleted once collide with player.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'When enemy is killed, some of them should drop pickable sushi.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'When player pick up sushi, they should be healed up as shown by the health bar. If they reach maximum health then nothing should happen.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'When enemy bullet shoot directly to player, it should be deleted once it hits the player character mesh.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'When player walk to enemy, it should not be damaged as long as the mesh is not colliding to enemy’s aabb.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'AI tests'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Stay far away from the enemy. The enemy should be moving and idling randomly.'
$p.Range.ListFormat.ListLevelNumber = 2
$p.LeftIndent = 72
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Stay closer to enemy but hide behind an object. The enemy should try to find the player with A* algorithm.'
$p.Range.ListFormat.ListLevelNumber = 2
$p.LeftIndent = 72
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Stay out of the object. The enemy should move closer to player and start shooting from a distance.'
$p.Range.ListFormat.ListLevelNumber = 2
$p.LeftIndent = 72
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'When stay still, player idling animation should be played.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Bees and bomber should face the direction of player.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Wolves should face to the direction they are travelling.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18

$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Text = 'Walk around an object (a pillar) and see the layers of rendering. Every entity should be drawn that the one in front is rendered afterwards.'
$p.Range.ListFormat.ListLevelNumber = 1
$p.LeftIndent = 36
$p.FirstLineIndent = -18
